$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 61

# Column A holds a date-looking value but must stay a plain text string
# (matching existing rows, which are stored as text, not real dates).
# A leading quote-prefix forces Excel to treat it as text instead of
# auto-converting it to a date serial number.
$ws.Cells.Item($row, 1).Value = "'01/24/2026"

$ws.Cells.Item($row, 2).Value = 11873.3
$ws.Cells.Item($row, 3).Value = 0.2328945415145907
$ws.Cells.Item($row, 4).Value = 0.7671054584854093
$ws.Cells.Item($row, 5).Value = -186.79
$ws.Cells.Item($row, 6).Value = -25.41
$ws.Cells.Item($row, 7).Value = -21683.68
$ws.Cells.Item($row, 8).Value = -70.42
$ws.Cells.Item($row, 9).Value = -437.62
$ws.Cells.Item($row, 10).Value = -13.66
$ws.Cells.Item($row, 11).Value = -22121.3
$ws.Cells.Item($row, 12).Value = -65.06999999999999
